$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix T2 and T7: numeric 0 -> text "0" (matches the rest of column T) ---
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "0"

$ws.Range("T7").NumberFormat = "@"
$ws.Range("T7").Value = "0"

# --- Row 8 (developer role): update Definition text ---
$ws.Range("D8").Value = "A role that involves having some responsibility for a development process."

# --- Row 11 (development process): update Definition text ---
$ws.Range("D11").Value = "A process that is creation, modification or tailored application of a project, service, intervention, commodity, product or enterprise."

# --- Delete row 13 ("patient and public involvement") entirely; rows below shift up ---
$ws.Rows("13:13").Delete()

# After the shift, old rows 14-20 are now rows 13-19. Fill in the blank ID/Definition
# cells (and a few other columns) for what are now rows 13-17.

# Row 13 (was row 14): patient and public involvement and engagement
$ws.Range("A13").Value = "BCIO:050277"
$ws.Range("D13").Value = "A development process in which patients or members of the public participate by virtue of their stakeholder role."

# Row 14 (was row 15): product development process
$ws.Range("A14").Value = "BCIO:050278"
$ws.Range("D14").Value = "A development process that is of a product."

# Row 15 (was row 16): project development process
$ws.Range("A15").Value = "BCIO:050279"
$ws.Range("D15").Value = "A development process that is of a project."

# Row 16 (was row 17): service development process
$ws.Range("A16").Value = "BCIO:050280"
$ws.Range("D16").Value = "A development process that is of a service."

# Row 17 (was row 18): stakeholder
$ws.Range("A17").Value = "BCIO:050276"
$ws.Range("D17").Value = "An agent that has a stakeholder role."
$ws.Range("G17").Value = "agent"
$ws.Range("H17").Value = "independent continuant"
$ws.Range("M17").Value = "stakeholder role"
